$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stimuli order: each row gives the image path (col B), the German verb
# (col C) and the category (col D) for data rows 2-33 (A holds the 0-based
# index and does not change).
$rows = @(
    @{ Row=2;  B="car/car078.png"; C="nehmen";    D="car" },
    @{ Row=3;  B="car/car076.png"; C="scheitern"; D="car" },
    @{ Row=4;  B="car/car087.png"; C="rasen";     D="car" },
    @{ Row=5;  B="car/car116.png"; C="runden";    D="car" },
    @{ Row=6;  B="car/car072.png"; C="kaufen";    D="car" },
    @{ Row=7;  B="car/car080.png"; C="währen";    D="car" },
    @{ Row=8;  B="car/car074.png"; C="liefern";   D="car" },
    @{ Row=9;  B="dog/dog091.png"; C="fliegen";   D="dog" },
    @{ Row=10; B="car/car101.png"; C="sieben";    D="car" },
    @{ Row=11; B="dog/dog116.png"; C="antun";     D="dog" },
    @{ Row=12; B="dog/dog115.png"; C="opfern";    D="dog" },
    @{ Row=13; B="car/car092.png"; C="hauen";     D="car" },
    @{ Row=14; B="dog/dog082.png"; C="krachen";   D="dog" },
    @{ Row=15; B="dog/dog107.png"; C="pflegen";   D="dog" },
    @{ Row=16; B="dog/dog067.png"; C="backen";    D="dog" },
    @{ Row=17; B="dog/dog106.png"; C="loben";     D="dog" },
    @{ Row=18; B="dog/dog066.png"; C="tagen";     D="dog" },
    @{ Row=19; B="dog/dog112.png"; C="formen";    D="dog" },
    @{ Row=20; B="dog/dog081.png"; C="ehren";     D="dog" },
    @{ Row=21; B="car/car089.png"; C="fesseln";   D="car" },
    @{ Row=22; B="dog/dog114.png"; C="klappen";   D="dog" },
    @{ Row=23; B="car/car086.png"; C="füllen";    D="car" },
    @{ Row=24; B="dog/dog095.png"; C="strahlen";  D="dog" },
    @{ Row=25; B="dog/dog125.png"; C="starten";   D="dog" },
    @{ Row=26; B="dog/dog104.png"; C="laufen";    D="dog" },
    @{ Row=27; B="car/car081.png"; C="töten";     D="car" },
    @{ Row=28; B="car/car071.png"; C="hoffen";    D="car" },
    @{ Row=29; B="car/car094.png"; C="stechen";   D="car" },
    @{ Row=30; B="dog/dog093.png"; C="raten";     D="dog" },
    @{ Row=31; B="car/car085.png"; C="rücken";    D="car" },
    @{ Row=32; B="dog/dog096.png"; C="biegen";    D="dog" },
    @{ Row=33; B="car/car102.png"; C="dauern";    D="car" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
